# Updated symbol list on Sun Jan 15 07:47:59 UTC 2023 with GitHub Actions
# Refreshes Price/Volume(1h) columns for existing rows and re-syncs the
# coin ranking block (rows 16-22) whose order shifted upstream.
# NumberFormat "@" (Text) is applied before each D/E write so the
# numeric-looking strings ("296.87", "-4.49%") stay stored as text
# (matching the sheet's existing inlineStr cells) instead of being
# auto-coerced into numbers/percentages by Excel; the format is reset to
# the sheet's default "Normal" style right after so no stray formatting
# is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value2 = '296.87'
$ws.Cells.Item(2,4).Style = "Normal"

$ws.Cells.Item(2,5).NumberFormat = "@"
$ws.Cells.Item(2,5).Value2 = '-4.49%'
$ws.Cells.Item(2,5).Style = "Normal"

$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value2 = '31.63'
$ws.Cells.Item(3,4).Style = "Normal"

$ws.Cells.Item(3,5).NumberFormat = "@"
$ws.Cells.Item(3,5).Value2 = '-0.69%'
$ws.Cells.Item(3,5).Style = "Normal"

$ws.Cells.Item(4,4).NumberFormat = "@"
$ws.Cells.Item(4,4).Value2 = '5.128'
$ws.Cells.Item(4,4).Style = "Normal"

$ws.Cells.Item(4,5).NumberFormat = "@"
$ws.Cells.Item(4,5).Value2 = '-4.21%'
$ws.Cells.Item(4,5).Style = "Normal"

$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value2 = '0.07486'
$ws.Cells.Item(5,4).Style = "Normal"

$ws.Cells.Item(5,5).NumberFormat = "@"
$ws.Cells.Item(5,5).Value2 = '-1.16%'
$ws.Cells.Item(5,5).Style = "Normal"

$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value2 = '7.741'
$ws.Cells.Item(6,4).Style = "Normal"

$ws.Cells.Item(6,5).NumberFormat = "@"
$ws.Cells.Item(6,5).Value2 = '-1.34%'
$ws.Cells.Item(6,5).Style = "Normal"

$ws.Cells.Item(7,4).NumberFormat = "@"
$ws.Cells.Item(7,4).Value2 = '1.714'
$ws.Cells.Item(7,4).Style = "Normal"

$ws.Cells.Item(7,5).NumberFormat = "@"
$ws.Cells.Item(7,5).Value2 = '5.52%'
$ws.Cells.Item(7,5).Style = "Normal"

$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value2 = '3.799'
$ws.Cells.Item(8,4).Style = "Normal"

$ws.Cells.Item(8,5).NumberFormat = "@"
$ws.Cells.Item(8,5).Value2 = '2.31%'
$ws.Cells.Item(8,5).Style = "Normal"

$ws.Cells.Item(9,4).NumberFormat = "@"
$ws.Cells.Item(9,4).Value2 = '0.9317'
$ws.Cells.Item(9,4).Style = "Normal"

$ws.Cells.Item(9,5).NumberFormat = "@"
$ws.Cells.Item(9,5).Value2 = '1.18%'
$ws.Cells.Item(9,5).Style = "Normal"

$ws.Cells.Item(10,4).NumberFormat = "@"
$ws.Cells.Item(10,4).Value2 = '0.1706'
$ws.Cells.Item(10,4).Style = "Normal"

$ws.Cells.Item(10,5).NumberFormat = "@"
$ws.Cells.Item(10,5).Value2 = '-0.66%'
$ws.Cells.Item(10,5).Style = "Normal"

$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value2 = '0.07165'
$ws.Cells.Item(11,4).Style = "Normal"

$ws.Cells.Item(11,5).NumberFormat = "@"
$ws.Cells.Item(11,5).Value2 = '-6.27%'
$ws.Cells.Item(11,5).Style = "Normal"

$ws.Cells.Item(12,4).NumberFormat = "@"
$ws.Cells.Item(12,4).Value2 = '0.07887'
$ws.Cells.Item(12,4).Style = "Normal"

$ws.Cells.Item(12,5).NumberFormat = "@"
$ws.Cells.Item(12,5).Value2 = '-3.96%'
$ws.Cells.Item(12,5).Style = "Normal"

$ws.Cells.Item(13,4).NumberFormat = "@"
$ws.Cells.Item(13,4).Value2 = '0.03013'
$ws.Cells.Item(13,4).Style = "Normal"

$ws.Cells.Item(13,5).NumberFormat = "@"
$ws.Cells.Item(13,5).Value2 = '-0.65%'
$ws.Cells.Item(13,5).Style = "Normal"

$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value2 = '0.09903'
$ws.Cells.Item(14,4).Style = "Normal"

$ws.Cells.Item(14,5).NumberFormat = "@"
$ws.Cells.Item(14,5).Value2 = '0.13%'
$ws.Cells.Item(14,5).Style = "Normal"

$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value2 = '0.001495'
$ws.Cells.Item(15,4).Style = "Normal"

$ws.Cells.Item(15,5).NumberFormat = "@"
$ws.Cells.Item(15,5).Value2 = '-2.82%'
$ws.Cells.Item(15,5).Style = "Normal"

$ws.Cells.Item(16,2).Value2 = 'CoinExToken'

$ws.Cells.Item(16,3).Value2 = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'

$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value2 = '0.04655'
$ws.Cells.Item(16,4).Style = "Normal"

$ws.Cells.Item(16,5).NumberFormat = "@"
$ws.Cells.Item(16,5).Value2 = '2.10%'
$ws.Cells.Item(16,5).Style = "Normal"

$ws.Cells.Item(17,2).Value2 = 'TigerCash'

$ws.Cells.Item(17,3).Value2 = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'

$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value2 = '0.006297'
$ws.Cells.Item(17,4).Style = "Normal"

$ws.Cells.Item(17,5).NumberFormat = "@"
$ws.Cells.Item(17,5).Value2 = '-4.39%'
$ws.Cells.Item(17,5).Style = "Normal"

$ws.Cells.Item(18,2).Value2 = 'LEO'

$ws.Cells.Item(18,3).Value2 = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'

$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value2 = '3.453'
$ws.Cells.Item(18,4).Style = "Normal"

$ws.Cells.Item(18,5).NumberFormat = "@"
$ws.Cells.Item(18,5).Value2 = '-1.19%'
$ws.Cells.Item(18,5).Style = "Normal"

$ws.Cells.Item(19,2).Value2 = 'BTSEToken'

$ws.Cells.Item(19,3).Value2 = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'

$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value2 = '2.228'
$ws.Cells.Item(19,4).Style = "Normal"

$ws.Cells.Item(19,5).NumberFormat = "@"
$ws.Cells.Item(19,5).Value2 = '-0.59%'
$ws.Cells.Item(19,5).Style = "Normal"

$ws.Cells.Item(20,2).Value2 = 'BitpandaEcosystemToken'

$ws.Cells.Item(20,3).Value2 = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'

$ws.Cells.Item(20,4).NumberFormat = "@"
$ws.Cells.Item(20,4).Value2 = '0.3274'
$ws.Cells.Item(20,4).Style = "Normal"

$ws.Cells.Item(20,5).NumberFormat = "@"
$ws.Cells.Item(20,5).Value2 = '-1.15%'
$ws.Cells.Item(20,5).Style = "Normal"

$ws.Cells.Item(21,2).Value2 = 'ProBitToken'

$ws.Cells.Item(21,3).Value2 = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'

$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value2 = '0.1327'
$ws.Cells.Item(21,4).Style = "Normal"

$ws.Cells.Item(21,5).NumberFormat = "@"
$ws.Cells.Item(21,5).Value2 = '-0.74%'
$ws.Cells.Item(21,5).Style = "Normal"

$ws.Cells.Item(22,2).Value2 = 'MCDex'

$ws.Cells.Item(22,3).Value2 = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'

$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value2 = '4.569'
$ws.Cells.Item(22,4).Style = "Normal"

$ws.Cells.Item(22,5).NumberFormat = "@"
$ws.Cells.Item(22,5).Value2 = '8.40%'
$ws.Cells.Item(22,5).Style = "Normal"

$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value2 = '0.1559'
$ws.Cells.Item(23,4).Style = "Normal"

$ws.Cells.Item(23,5).NumberFormat = "@"
$ws.Cells.Item(23,5).Value2 = '-4.32%'
$ws.Cells.Item(23,5).Style = "Normal"

$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value2 = '0.001219'
$ws.Cells.Item(24,4).Style = "Normal"

$ws.Cells.Item(24,5).NumberFormat = "@"
$ws.Cells.Item(24,5).Value2 = '-0.80%'
$ws.Cells.Item(24,5).Style = "Normal"

$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value2 = '0.004426'
$ws.Cells.Item(25,4).Style = "Normal"

$ws.Cells.Item(25,5).NumberFormat = "@"
$ws.Cells.Item(25,5).Value2 = '-1.38%'
$ws.Cells.Item(25,5).Style = "Normal"

$ws.Cells.Item(26,5).NumberFormat = "@"
$ws.Cells.Item(26,5).Value2 = '0.00%'
$ws.Cells.Item(26,5).Style = "Normal"

$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value2 = '0.0001877'
$ws.Cells.Item(27,4).Style = "Normal"

$ws.Cells.Item(27,5).NumberFormat = "@"
$ws.Cells.Item(27,5).Value2 = '7.76%'
$ws.Cells.Item(27,5).Style = "Normal"

$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value2 = '0.01681'
$ws.Cells.Item(39,4).Style = "Normal"

$ws.Cells.Item(39,5).NumberFormat = "@"
$ws.Cells.Item(39,5).Value2 = '-0.07%'
$ws.Cells.Item(39,5).Style = "Normal"

$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value2 = '0.04465'
$ws.Cells.Item(40,4).Style = "Normal"

$ws.Cells.Item(40,5).NumberFormat = "@"
$ws.Cells.Item(40,5).Value2 = '-2.44%'
$ws.Cells.Item(40,5).Style = "Normal"

$ws.Cells.Item(41,4).NumberFormat = "@"
$ws.Cells.Item(41,4).Value2 = '0.007125'
$ws.Cells.Item(41,4).Style = "Normal"

$ws.Cells.Item(41,5).NumberFormat = "@"
$ws.Cells.Item(41,5).Value2 = '-1.13%'
$ws.Cells.Item(41,5).Style = "Normal"

$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value2 = '0.1329'
$ws.Cells.Item(42,4).Style = "Normal"

$ws.Cells.Item(42,5).NumberFormat = "@"
$ws.Cells.Item(42,5).Value2 = '-2.66%'
$ws.Cells.Item(42,5).Style = "Normal"

$ws.Cells.Item(43,5).NumberFormat = "@"
$ws.Cells.Item(43,5).Value2 = '-8.85%'
$ws.Cells.Item(43,5).Style = "Normal"

$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value2 = '0.01132'
$ws.Cells.Item(44,4).Style = "Normal"

$ws.Cells.Item(44,5).NumberFormat = "@"
$ws.Cells.Item(44,5).Value2 = '-19.97%'
$ws.Cells.Item(44,5).Style = "Normal"

$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value2 = '0.00006038'
$ws.Cells.Item(45,4).Style = "Normal"

$ws.Cells.Item(45,5).NumberFormat = "@"
$ws.Cells.Item(45,5).Value2 = '-1.96%'
$ws.Cells.Item(45,5).Style = "Normal"

$ws.Cells.Item(46,4).NumberFormat = "@"
$ws.Cells.Item(46,4).Value2 = '1.918'
$ws.Cells.Item(46,4).Style = "Normal"

$ws.Cells.Item(46,5).NumberFormat = "@"
$ws.Cells.Item(46,5).Value2 = '1.34%'
$ws.Cells.Item(46,5).Style = "Normal"

$ws.Cells.Item(47,5).NumberFormat = "@"
$ws.Cells.Item(47,5).Value2 = '-0.30%'
$ws.Cells.Item(47,5).Style = "Normal"
